$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round ConvexHullArea values (column D, rows 2-27) to the nearest integer
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value2 = [Math]::Round([double]$cell.Value2, 0)
}
